$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New fake-data rows 2-7 (gender, first_name, Last_name, theDay, TheMonth, theYear, email, company, testPassword)
$data = @(
    @("male", "Joannie",    "Koelpin",  "30", "4", "2019", "gordon.grady@hotmail.com",  "Gottlieb Inc",               "zzklx9z4e"),
    @("male", "Miranda",    "Steuber",  "29", "4", "2019", "luther.goldner@gmail.com",  "Shields, Adams and Kemmer",  "zfjcn4g24jpf"),
    @("male", "Margarito",  "Corwin",   "29", "4", "2019", "cecil.smitham@yahoo.com",   "Ullrich, Bins and Sauer",    "4jzty0p6u"),
    @("male", "Antony",     "Schmeler", "7",  "5", "2019", "fred.sanford@gmail.com",    "Kihn, Gibson and Cremin",    "8bgpvktgwzop"),
    @("male", "Nicholas",   "Veum",     "29", "4", "2019", "lazaro.carter@yahoo.com",   "Cronin, King and Windler",   "ruf0s21109ffa"),
    @("male", "Dolly",      "Muller",   "3",  "5", "2019", "cyril.lueilwitz@gmail.com", "Rempel and Sons",            "dwdtwivdq2xnn")
)

# The theDay/TheMonth/theYear columns (D, E, F) hold numeric-looking text
# (e.g. "30", "4", "2019"); force them to be stored as text, not numbers.
$ws.Range("D2:F7").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}
